# Atualização automática de preços de eletricidade
# Updates row 2 (the single data row) of the spot price table with the
# latest values pulled from the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (date serial advances by one day)
$ws.Range("A2").Value = 45937

# Hourly prices
$ws.Range("B2").Value = 110.34
$ws.Range("C2").Value = 108.68
$ws.Range("D2").Value = 107.25
$ws.Range("E2").Value = 103.66
$ws.Range("F2").Value = 103.35
$ws.Range("G2").Value = 102.93
$ws.Range("H2").Value = 108.51
$ws.Range("I2").Value = 121.01
$ws.Range("J2").Value = 132.32
$ws.Range("K2").Value = 118.1
$ws.Range("L2").Value = 75.89
$ws.Range("M2").Value = 51.76
$ws.Range("N2").Value = 50.88
$ws.Range("O2").Value = 41.41
$ws.Range("P2").Value = 37.34
$ws.Range("Q2").Value = 42.81
$ws.Range("R2").Value = 54.95
$ws.Range("S2").Value = 77.72
$ws.Range("T2").Value = 109.65
$ws.Range("U2").Value = 150
$ws.Range("V2").Value = 168.35
$ws.Range("W2").Value = 156.29
$ws.Range("X2").Value = 119.76
$ws.Range("Y2").Value = 114.79
$ws.Range("Z2").Value = 98.66

# Slot summaries
$ws.Range("AB2").Value = 139.8
$ws.Range("AD2").Value = 162.32
$ws.Range("AF2").Value = 129.82
$ws.Range("AG2").Value = "10h-17h"
